# Apply updated cryptocurrency price/volume data (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.948.32'
$ws.Range('E2').Value = '  +2.68%  '
$ws.Range('D3').Value = '1.599.40'
$ws.Range('E3').Value = '  +2.73%  '
$ws.Range('E4').Value = '  -0.24%  '
$ws.Range('D5').Value = "'211.38"
$ws.Range('E5').Value = '  +2.30%  '
$ws.Range('E6').Value = '  -0.25%  '
$ws.Range('D7').Value = "'0.483"
$ws.Range('E7').Value = '  +1.23%  '
$ws.Range('E8').Value = '  +1.70%  '
$ws.Range('D9').Value = "'0.0613"
$ws.Range('E9').Value = '  +0.60%  '
$ws.Range('D10').Value = "'18.12"
$ws.Range('E10').Value = '  +1.69%  '
$ws.Range('D11').Value = "'0.0813"
$ws.Range('E11').Value = '  +3.96%  '
$ws.Range('D12').Value = '1.822.29'
$ws.Range('E12').Value = '  +2.77%  '
$ws.Range('D13').Value = '1.595.46'
$ws.Range('E13').Value = '  +2.47%  '
$ws.Range('D14').Value = "'4.01"
$ws.Range('E14').Value = '  +0.17%  '
$ws.Range('D15').Value = "'0.509"
$ws.Range('E15').Value = '  +0.80%  '
$ws.Range('D16').Value = '25.963.46'
$ws.Range('E16').Value = '  +2.72%  '
$ws.Range('D17').Value = "'60.20"
$ws.Range('E17').Value = '  +2.08%  '
$ws.Range('D18').Value = '0.0₃0721'
$ws.Range('E18').Value = '  +1.55%  '
$ws.Range('E19').Value = '  -0.22%  '
$ws.Range('D20').Value = "'199.99"
$ws.Range('E20').Value = '  +7.48%  '
$ws.Range('D21').Value = "'4.22"
$ws.Range('E21').Value = '  +2.21%  '
$ws.Range('E22').Value = '  -0.39%  '
$ws.Range('D23').Value = "'6.00"
$ws.Range('E23').Value = '  +2.52%  '
$ws.Range('E24').Value = '  +9.18%  '
$ws.Range('D25').Value = "'141.44"
$ws.Range('E25').Value = '  +0.97%  '
$ws.Range('E26').Value = '  -0.19%  '
$ws.Range('D27').Value = "'0.122"
$ws.Range('E27').Value = '  -5.26%  '
$ws.Range('E28').Value = '  +1.72%  '
$ws.Range('D29').Value = "'6.43"
$ws.Range('E29').Value = '  +0.43%  '
$ws.Range('E30').Value = '  +1.57%  '
$ws.Range('D31').Value = "'0.0474"
$ws.Range('E31').Value = '  +1.59%  '
$ws.Range('D32').Value = "'3.12"
$ws.Range('E32').Value = '  +2.53%  '
$ws.Range('D33').Value = "'2.96"
$ws.Range('E33').Value = '  -0.41%  '
$ws.Range('E34').Value = '  +1.04%  '
$ws.Range('E35').Value = '  +1.03%  '
$ws.Range('B36').Value = 'Maker'
$ws.Range('C36').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D36').Value = '1.128.42'
$ws.Range('E36').Value = '  +4.47%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = "'0.0165"
$ws.Range('E37').Value = '  +11.16%  '
$ws.Range('E38').Value = '  +0.15%  '
$ws.Range('D39').Value = "'0.792"
$ws.Range('E39').Value = '  +4.26%  '
$ws.Range('D40').Value = "'2.31"
$ws.Range('E40').Value = '  +0.40%  '
$ws.Range('D41').Value = "'0.490"
$ws.Range('E41').Value = '  -0.31%  '
$ws.Range('D42').Value = "'0.783"
$ws.Range('E42').Value = '  -1.49%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = "'5.16"
$ws.Range('E43').Value = '  +1.95%  '
$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').Value = '1.733.69'
$ws.Range('E44').Value = '  +2.69%  '
$ws.Range('D45').Value = "'92.85"
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('D46').Value = "'1.52"
$ws.Range('E46').Value = '  +3.70%  '
$ws.Range('D47').Value = "'53.29"
$ws.Range('E47').Value = '  +1.71%  '
$ws.Range('D48').Value = "'0.0503"
$ws.Range('E48').Value = '  +0.03%  '
$ws.Range('E49').Value = '  +0.99%  '
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('D51').Value = '0.0₇0922'
$ws.Range('E51').Value = '  -17.34%  '
